$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 210 (shifts rows 210:278 down to 211:279).
$ws.Rows("210:210").Insert()

# The freshly inserted row inherited mangled formatting from the insert
# operation; restore the exact formats used by the surrounding "content"
# rows (e.g. row 215, which has the same style signature as the target
# row) by copying formats only - this keeps cell styles 7/4/29/19/20/35
# instead of new/duplicated style indices.
$ws.Range("A215:H215").Copy()
$ws.Range("A210:H210").PasteSpecial(-4122)

# Populate the new row with the new API entry.
$ws.Cells.Item(210, 2).Value = "transaction.synchronize.project.setProjectSection"
$ws.Cells.Item(210, 3).Value = "Menyinkronkan Data Project Section"

# Clear clipboard marquee.
$excel.CutCopyMode = $false

# Restore the selected cell in the (now-shifted) view to C201, matching
# the post-edit selection recorded for the sheet.
$ws.Range("C201").Select()
